$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing column E (password) into the new column F first,
# so the new "fakultas" column can be inserted at E without disturbing
# formatting (a real Columns.Insert() would copy D's hyperlink style).
for ($r = 1; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 5).Value()
}

# Column E header + fakultas values for each Dosen row.
$fakultas = @(
    "fakultas",
    "Pendidikan Profesi Psikologi",
    "`tProfesi Akuntan",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika",
    "`tSistem Informasi",
    "`tTeknik Informatika"
)

for ($r = 1; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = $fakultas[$r - 1]
}

# Match the saved selection from the authored workbook.
$ws.Range("E12").Select()
